# Updated vehicle data to include gearbox parameters
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 28 (before the blank separator row 29),
# shifting everything from the old row 29 down to row 31.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# gear_ratio value changed from 11 to 12, and its "Value" style flips from Good to Neutral
$ws.Range("B24").Value = 12
$ws.Range("B24").Style = "Neutral"

# Row 25: motor_inertia -> j_m (value/style/notes text unchanged)
$ws.Range("A25").Value = "j_m"
$ws.Range("B25").Value = 0.000274
$ws.Range("B25").Style = "Good"
$ws.Range("C25").Value = "electric motor rotational inertia"

# Row 26 (new): J_s1, input shaft inertia of gearbox (zero, motor is the input shaft)
$ws.Range("A26").Value = "J_s1"
$ws.Range("B26").Value = 0
$ws.Range("B26").Style = "Bad"
$ws.Range("C26").Value = "inertia of gearbox input shaft (motor is the input shaft so this is zero)"

# Row 27: gearbox_inertia -> J_1
$ws.Range("A27").Value = "J_1"
$ws.Range("B27").Value = 0.0008
$ws.Range("B27").Style = "Bad"
$ws.Range("C27").Value = "inertia of gearbox components as seen by the input shaft"

# Row 28 (new): J_2, gearbox output shaft inertia (zero, using J_1 instead)
$ws.Range("A28").Value = "J_2"
$ws.Range("B28").Value = 0
$ws.Range("B28").Style = "Bad"
$ws.Range("C28").Value = "inertia of gearbox components as seen by the output shaft (using J_1 so this must be zero)"

# Row 29: output_shaft_inertia -> J_s2
$ws.Range("A29").Value = "J_s2"
$ws.Range("B29").Value = 0.0005
$ws.Range("B29").Style = "Bad"
$ws.Range("C29").Value = "inertia of ouput shaft components not included in gearbox_inertia (hub, brake rotor, lug nuts, etc)"

# Row 30: wheel_inertia -> J_w
$ws.Range("A30").Value = "J_w"
$ws.Range("B30").Value = 0.003
$ws.Range("B30").Style = "Bad"
$ws.Range("C30").Value = "inertia of rim and tire"

# Update view settings
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D7").Select()

$wb.Save()
